$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").NumberFormat = "General"
$ws.Range("C2").Value = 5118275524
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "General"
$ws.Range("C3").Value = 5118275524
$ws.Range("C3").NumberFormat = "@"
